# Apply refreshed capital-structure figures for the France Brokerage & Investment Banking rows (rows 2-4).
# Row 3 (Bourse Direct) and row 4 (VIEL & Cie) swap position, and every metric column is refreshed
# with newly pulled figures; a few trailing columns (AN/AO/AP/AQ) are cleared or added per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0244
$ws.Range("E2").Value = 0.1585
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 68.32000000000001
$ws.Range("L2").Value = 0.06338250301512201
$ws.Range("M2").Value = 18.753
$ws.Range("N2").Value = 0.02760229614365617
$ws.Range("O2").Value = 0.2744877049180328
$ws.Range("P2").Value = 18.5
$ws.Range("Q2").Value = 0.02722990874300854
$ws.Range("R2").Value = 0.2707845433255269
$ws.Range("S2").Value = 0.2530000000000001
$ws.Range("T2").Value = 0.01349117474537408
$ws.Range("U2").Value = 462.4
$ws.Range("V2").Value = 0.6806005298793053
$ws.Range("W2").Value = 0.1166236779489123
$ws.Range("X2").Value = 0.04554896790699149
$ws.Range("Y2").Value = 0.07107471004192077
$ws.Range("Z2").Value = 1.658818097876269
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.0317500808745322
$ws.Range("AC2").Value = -0.0317500808745322
$ws.Range("AD2").Value = 707.5
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 707.5
$ws.Range("AG2").Value = 245.1
$ws.Range("AH2").Value = 0.5101305068858605
$ws.Range("AI2").Value = 0.5291697830964847
$ws.Range("AJ2").Value = 0.2651162790697675
$ws.Range("AK2").Value = 0.2802423965241254
$ws.Range("AM2").Value = -2.78
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# Row 3
$ws.Range("B3").Value = "Bourse Direct SA (ENXTPA:BSD)"
$ws.Range("D3").Value = 0.014
$ws.Range("E3").Value = 0.146
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5.62
$ws.Range("L3").Value = 0.1260089686098655
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0.08017118402282454
$ws.Range("X3").Value = 0.04138898506027669
$ws.Range("Y3").Value = 0.03878219896254785
$ws.Range("Z3").Value = 0.210576015108593
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.03163314704728253
$ws.Range("AC3").Value = -0.03163314704728253
$ws.Range("AD3").Value = 133.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 133.5
$ws.Range("AG3").Value = 133.5
$ws.Range("AH3").Value = 0.3997005988023952
$ws.Range("AI3").Value = 0.6443050193050194
$ws.Range("AJ3").Value = 0.3997005988023952
$ws.Range("AK3").Value = 0.6443050193050194
$ws.Range("AM3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AQ3").ClearContents()

# Row 4
$ws.Range("B4").Value = "VIEL & Cie, société anonyme (ENXTPA:VIL)"
$ws.Range("D4").Value = 0.0348
$ws.Range("E4").Value = 0.171
$ws.Range("K4").Value = 62.7
$ws.Range("L4").Value = 0.06067937675408885
$ws.Range("M4").Value = 18.753
$ws.Range("N4").Value = 0.03915848820212988
$ws.Range("O4").Value = 0.2990909090909091
$ws.Range("P4").Value = 18.5
$ws.Range("Q4").Value = 0.03863019419503028
$ws.Range("R4").Value = 0.2950558213716108
$ws.Range("S4").Value = 0.2530000000000001
$ws.Range("T4").Value = 0.01349117474537408
$ws.Range("U4").Value = 462.4
$ws.Range("V4").Value = 0.9655460430152433
$ws.Range("W4").Value = 0.153076171875
$ws.Range("X4").Value = 0.04970895075370629
$ws.Range("Y4").Value = 0.1033672211212937
$ws.Range("Z4").Value = 2.359132420091324
$ws.Range("AB4").Value = 0.03186701470178187
$ws.Range("AC4").Value = -0.03186701470178187
$ws.Range("AD4").Value = 574
$ws.Range("AF4").Value = 574
$ws.Range("AG4").Value = 111.6
$ws.Range("AH4").Value = 0.545160983949093
$ws.Range("AI4").Value = 0.5080545229244114
$ws.Range("AJ4").Value = 0.1889923793395428
$ws.Range("AK4").Value = 0.1672160623314355
$ws.Range("AM4").Value = -2.78
$ws.Range("AQ4").Value = -0
